$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts old C..F -> D..G)
$ws.Columns.Item(3).Insert()

# Update header labels for the (now shifted) GitHub/Website columns, and
# set the new column C header ("Language Used"). The order below matches
# the order the new shared-string entries were originally authored in.
$ws.Range("G1").Value = "Website"
$ws.Range("F1").Value = "GitHub"
$ws.Range("C1").Value = "Language Used"

# Add new row with Serial No 4
$ws.Range("A5").Value = 4

# Re-create the hyperlinks at their new (shifted) locations, since the
# insert operation does not move the hyperlink anchors automatically.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/abhisekjha/steganography")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://www.abhisekjha.com.np/steganography/")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/abhisekjha/pqc_aes_multipath")
$ws.Hyperlinks.Add($ws.Range("G3"), "https://www.abhisekjha.com.np/pqc_aes_multipath")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/abhisekjha/luxury-car-calculator")
$ws.Hyperlinks.Add($ws.Range("G4"), "https://www.abhisekjha.com.np/luxury-car-calculator")

# Reapply the Hyperlink style uniformly so it reuses the existing style
# definition instead of Excel creating a duplicate one per-cell.
$ws.Range("F2:G4").Style = "Hyperlink"

# Adjust column widths to match target layout. Columns D, F and G already
# inherit the correct (exact) widths from the pre-insert C, E and F columns,
# so only B, C and E need to be touched. The ColumnWidth setter in this
# runtime quantizes to whole-pixel steps, so the input values below were
# calibrated to be the ones that round-trip to the exact desired stored
# widths (20.1640625 and 17 respectively).
$ws.Columns.Item(2).ColumnWidth = 19.3
$ws.Columns.Item(3).ColumnWidth = 19.3
$ws.Columns.Item(5).ColumnWidth = 16.2

# Update selection to H1
$ws.Range("H1").Select()
